$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New helper block (rows 22-23): RGB components of the accent color used
# on the chart (186, 12, 47 == theme Accent1 == #BA0C2F) plus their
# 0-1 fractions for reuse elsewhere.
$ws.Range("O22").Interior.ThemeColor = 5

$ws.Range("Q22").Value = 186
$ws.Range("R22").Value = 12
$ws.Range("S22").Value = 47

$ws.Range("Q23").Formula = "=Q22/255"
$ws.Range("R23").Formula = "=R22/255"
$ws.Range("S23").Formula = "=S22/255"

# Move the active selection/scroll position to reflect the new data below
# the original table (matches the saved view state after the edit).
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("Q24").Select()
